$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Init")

$ws.Range("D5").Value = "A31"
$ws.Range("D6").Value = "B31"
$ws.Range("D7").Value = "C31"
$ws.Range("D8").Value = "G31"
$ws.Range("D9").Value = "H31"
$ws.Range("D10").Value = "I31"
$ws.Range("D11").Value = "J31"

$ws.Range("D15").Value = "E258"

$ws.Range("B20").Value = "cs3rpt2022_all_demand_units_v20251130.xlsx"

$ws.Range("D22").Value = "O440"

$ws.Range("E5:E11").Font.Bold = $true
$ws.Range("E5:E11").Font.Italic = $true

$ws.Range("D5:D11").Select()
$excel.ActiveWindow.ScrollRow = 3
